$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'307.39"
$ws.Cells.Item(2,4).Style = "Normal"

$ws.Cells.Item(3,4).Value = "'38.58"
$ws.Cells.Item(3,5).Value = "'8.28%"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Style = "Normal"

$ws.Cells.Item(4,2).Value = "HuobiToken"
$ws.Cells.Item(4,3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(4,4).Value = "'5.103"
$ws.Cells.Item(4,5).Value = "'1.20%"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Style = "Normal"

$ws.Cells.Item(5,2).Value = "Cronos"
$ws.Cells.Item(5,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(5,4).Value = "'0.08113"
$ws.Cells.Item(5,5).Value = "'1.30%"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Style = "Normal"

$ws.Cells.Item(6,2).Value = "FTXToken"
$ws.Cells.Item(6,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(6,4).Value = "'1.958"
$ws.Cells.Item(6,5).Value = "'4.98%"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Style = "Normal"

$ws.Cells.Item(7,2).Value = "KuCoinToken"
$ws.Cells.Item(7,3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(7,4).Value = "'7.955"
$ws.Cells.Item(7,5).Value = "'2.20%"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Style = "Normal"

$ws.Cells.Item(8,2).Value = "MXToken"
$ws.Cells.Item(8,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8,4).Value = "'0.9306"
$ws.Cells.Item(8,5).Value = "'1.02%"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Style = "Normal"

$ws.Cells.Item(9,2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9,4).Value = "'0.1409"
$ws.Cells.Item(9,5).Value = "'9.88%"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Style = "Normal"

$ws.Cells.Item(10,2).Value = "WazirX"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10,4).Value = "'0.1956"
$ws.Cells.Item(10,5).Value = "'2.75%"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Style = "Normal"

$ws.Cells.Item(11,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11,4).Value = "'0.09046"
$ws.Cells.Item(11,5).Value = "'-1.25%"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Style = "Normal"

$ws.Cells.Item(12,2).Value = "BitrueCoin"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12,4).Value = "'0.03506"
$ws.Cells.Item(12,5).Value = "'2.83%"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Style = "Normal"

$ws.Cells.Item(13,2).Value = "BitMartToken"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13,4).Value = "'0.09827"
$ws.Cells.Item(13,5).Value = "'-0.23%"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Style = "Normal"

$ws.Cells.Item(14,2).Value = "BitForexToken"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14,4).Value = "'0.001405"
$ws.Cells.Item(14,5).Value = "'0.00%"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Style = "Normal"

$ws.Cells.Item(15,2).Value = "TigerCash"
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15,4).Value = "'0.006151"
$ws.Cells.Item(15,5).Value = "'-1.01%"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Style = "Normal"

$ws.Cells.Item(16,2).Value = "LEO"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16,4).Value = "'3.729"
$ws.Cells.Item(16,5).Value = "'-3.09%"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Style = "Normal"

$ws.Cells.Item(17,2).Value = "GateToken"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17,4).Value = "'4.198"
$ws.Cells.Item(17,5).Value = "'1.50%"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Style = "Normal"

$ws.Cells.Item(18,4).Value = "'3.419"
$ws.Cells.Item(18,5).Value = "'5.23%"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Style = "Normal"

$ws.Cells.Item(19,4).Value = "'0.3463"
$ws.Cells.Item(19,5).Value = "'1.26%"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Style = "Normal"

$ws.Cells.Item(20,4).Value = "'0.1342"
$ws.Cells.Item(20,5).Value = "'-0.47%"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Style = "Normal"

$ws.Cells.Item(21,4).Value = "'4.797"
$ws.Cells.Item(21,5).Value = "'-7.29%"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Style = "Normal"

$ws.Cells.Item(22,5).Value = "'6.32%"
$ws.Cells.Item(22,5).Style = "Normal"

$ws.Cells.Item(23,4).Value = "'0.04415"
$ws.Cells.Item(23,5).Value = "'-0.29%"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Style = "Normal"

$ws.Cells.Item(24,4).Value = "'0.001221"
$ws.Cells.Item(24,5).Value = "'-1.02%"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Style = "Normal"

$ws.Cells.Item(25,5).Value = "'-0.96%"
$ws.Cells.Item(25,5).Style = "Normal"

$ws.Cells.Item(27,4).Value = "'0.0001302"
$ws.Cells.Item(27,5).Value = "'4.00%"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Style = "Normal"

$ws.Cells.Item(39,4).Value = "'0.02077"
$ws.Cells.Item(39,5).Value = "'7.24%"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Style = "Normal"

$ws.Cells.Item(40,4).Value = "'0.05130"
$ws.Cells.Item(40,5).Value = "'-1.93%"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Style = "Normal"

$ws.Cells.Item(41,4).Value = "'0.007477"
$ws.Cells.Item(41,5).Value = "'-1.77%"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Style = "Normal"

$ws.Cells.Item(43,4).Value = "'0.1356"
$ws.Cells.Item(43,5).Value = "'0.40%"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Style = "Normal"

$ws.Cells.Item(44,4).Value = "'0.002132"
$ws.Cells.Item(44,5).Value = "'-1.39%"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Style = "Normal"

$ws.Cells.Item(45,4).Value = "'0.009272"
$ws.Cells.Item(45,5).Value = "'-3.73%"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Style = "Normal"

$ws.Cells.Item(46,4).Value = "'0.00006253"
$ws.Cells.Item(46,5).Value = "'0.90%"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Style = "Normal"

$ws.Cells.Item(47,4).Value = "'0.00000000751"
$ws.Cells.Item(47,5).Value = "'0.05%"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Style = "Normal"

$ws.Cells.Item(48,4).Value = "'0.003036"
$ws.Cells.Item(48,4).Style = "Normal"

$ws.Cells.Item(49,5).Value = "'-3.51%"
$ws.Cells.Item(49,5).Style = "Normal"

$ws.Cells.Item(50,5).Value = "'0.05%"
$ws.Cells.Item(50,5).Style = "Normal"

$ws.Cells.Item(51,5).Value = "'0.05%"
$ws.Cells.Item(51,5).Style = "Normal"
